{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two paragraphs that need their runs (and the spell-check\n// proofErr wrapped run in the middle) collapsed into a single plain run.\nlet seedParagraph = null;\nlet examplesParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text;\n  if (seedParagraph === null && text.indexOf(\"If seed is None\") === 0) {\n    seedParagraph = p;\n  }\n  if (examplesParagraph === null && text.indexOf(\"Examples provided\") === 0) {\n    examplesParagraph = p;\n  }\n}\n\nif (!seedParagraph || !examplesParagraph) {\n  throw new Error(\"Could not locate expected paragraphs in document body.\");\n}\n\n// Collapse \"If seed is None in any method, it will revert to the \" +\n// \"Genertor's\" (proofErr-wrapped) + \" default which is 0.\" into one run.\nseedParagraph\n  .getRange()\n  .insertText(\n    \"If seed is None in any method, it will revert to the Genertor\\u2019s default which is 0.\",\n    \"Replace\"\n  );\n\n// Collapse \"Examples provided for console, \" + \"Jupyter\" (proofErr-wrapped) +\n// \" Notebooks, and GUI.\" into one run.\nexamplesParagraph\n  .getRange()\n  .insertText(\"Examples provided for console, Jupyter Notebooks, and GUI.\", \"Replace\");\n\nawait context.sync();\n\n// Append two new bulleted list items after the \"Examples provided...\" entry:\n//   * LED:\n//     o I moved reading the CSV inside the method.\nconst ledParagraph = examplesParagraph.insertParagraph(\"LED:\", \"After\");\nconst csvParagraph = ledParagraph.insertParagraph(\n  \"I moved reading the CSV inside the method.\",\n  \"After\"\n);\ncsvParagraph.listItem.level = 1;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-ParaIndexStartingWith($doc, $prefix) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        if ($t.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# --- Collapse \"If seed is None ... Genertor's default which is 0.\" (which\n# was split across 3 runs because of a spell-check proofErr wrap around\n# \"Genertor's\") into a single plain run. ---\n$seedIdx = Get-ParaIndexStartingWith $d \"If seed is None\"\n$seedPara = $d.Paragraphs.Item($seedIdx)\n$seedRange = $seedPara.Range\n$seedRange.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark\n[void]$seedRange.Delete()\n$seedRange.InsertAfter(\"If seed is None in any method, it will revert to the Genertor\" + [char]0x2019 + \"s default which is 0.\") | Out-Null\n\n# --- Collapse \"Examples provided for console, Jupyter Notebooks, and GUI.\"\n# (split across 3 runs because of a spell-check proofErr wrap around\n# \"Jupyter\") into a single plain run. ---\n$examplesIdx = Get-ParaIndexStartingWith $d \"Examples provided\"\n$examplesPara = $d.Paragraphs.Item($examplesIdx)\n$examplesRange = $examplesPara.Range\n$examplesRange.MoveEnd(1, -1) | Out-Null\n[void]$examplesRange.Delete()\n$examplesRange.InsertAfter(\"Examples provided for console, Jupyter Notebooks, and GUI.\") | Out-Null\n\n# --- Add two new bulleted list items after it:\n#       * LED:\n#         o I moved reading the CSV inside the method.\n$examplesPara.Range.InsertParagraphAfter() | Out-Null\n\n$ledIdx = $examplesIdx + 1\n$ledPara = $d.Paragraphs.Item($ledIdx)\n$ledRange = $ledPara.Range\n$ledRange.MoveEnd(1, -1) | Out-Null\n$ledRange.Text = \"LED:\"\n\n$ledPara.Range.InsertParagraphAfter() | Out-Null\n\n$csvIdx = $ledIdx + 1\n$csvPara = $d.Paragraphs.Item($csvIdx)\n$csvRange = $csvPara.Range\n$csvRange.MoveEnd(1, -1) | Out-Null\n$csvRange.Text = \"I moved reading the CSV inside the method.\"\n$csvPara.Range.ListFormat.ListLevelNumber = 2   # 1-based => w:ilvl=1 (sub-bullet)\n"}
